# Second commit of the Project
# Updates a handful of test-data values on the "testDataAPI" sheet:
#  - refresh the email addresses used by three customer-creation rows
#  - refresh the email + generated customer id used by the delete-customer row
#  - move the active selection from A20 (now empty/out of range) to A13

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testDataAPI")

# Row 3 (Raunak111) - email refreshed
$ws.Range("B3").Value = "rao112@zmail.com"

# Row 4 (ashutosh111) - email refreshed
$ws.Range("B4").Value = "ahpr112@zmail.com"

# Row 5 (imraan111) - email refreshed
$ws.Range("B5").Value = "imr112@zmail.com"

# Row 9 (Atul) - email refreshed
$ws.Range("B9").Value = "at12@zmail.com"

# Row 13 - generated Stripe customer id refreshed
$ws.Range("A13").Value = "cus_JQcHyOgx0NzdKA"

# Move the sheet's active selection to A13
$ws.Activate()
$ws.Range("A13").Select()
